$wb = $excel.ActiveWorkbook

# --- Update the "Conversión del día" note on Hoja1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$newLine = [char]10
$newText = "Conversión del día 💰" + $newLine + `
    "✅ Dólar paralelo: 68" + $newLine + `
    $newLine + `
    "Binance" + $newLine + `
    "✅ 1000 Bs = 6.62 = 26179.23 pesos" + $newLine + `
    "✅ 26179.23 pesos = 6.61 = 979.33 Bs" + $newLine + `
    $newLine + `
    "Promedio competencia" + $newLine + `
    "✅ Tasa pesos: 20" + $newLine + `
    "✅ Tasa Bs: 20" + $newLine + `
    "✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $newText

# --- Update the rate figures on the "tasas" sheet ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 151.15
$ws2.Range("O10").Value = 3956.99
$ws2.Range("N12").Value = 3959
$ws2.Range("O12").Value = 148.101
